$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Isolate the run that currently reads:
#    " skills. Gained through my extensive experience performing in a
#      campanology (Church Bell Ringer) group"
#    from its neighbours ("Teamwork and group cooperation" before it and
#    " where we" after it) by dropping temporary bookmarks at both edges.
#    Word always breaks a run where a bookmark is inserted, and - important -
#    that split survives even after the temporary bookmark is deleted again.
# ---------------------------------------------------------------------------

$anchorBefore = "Teamwork and group cooperation"
$oldPhrase    = "Church Bell Ringer"
$anchorAfterOld = $oldPhrase + ") group"

$full = $d.Content.Text
$beforeEnd = $full.IndexOf($anchorBefore) + $anchorBefore.Length
$rBeforeEnd = $d.Range($beforeEnd, $beforeEnd)
$d.Bookmarks.Add("ZZZ_splitBefore", $rBeforeEnd) | Out-Null

$full = $d.Content.Text
$afterEnd = $full.IndexOf($anchorAfterOld) + $anchorAfterOld.Length
$rAfterEnd = $d.Range($afterEnd, $afterEnd)
$d.Bookmarks.Add("ZZZ_splitAfter", $rAfterEnd) | Out-Null

# ---------------------------------------------------------------------------
# 2. Replace "Church Bell Ringer" with "tower bell ringing" inside the now
#    fully isolated run, so the edit cannot bleed into neighbouring runs.
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$oldIdx = $full.IndexOf($oldPhrase)
$rOld = $d.Range($oldIdx, $oldIdx + $oldPhrase.Length)
$rOld.Text = "tower bell ringing"

# ---------------------------------------------------------------------------
# 3. Re-split the isolated run into the final fragments:
#      " skills. Gained through"
#      " my extensive experience performing in a campanology ("
#      "tower"
#      " "
#      "b"
#      "ell "
#      "r"
#      "ing"
#      "ing"
#      ") group"
#    A real "_GoBack" bookmark is dropped exactly where the user's cursor
#    was left (right after "through"); this also relocates the single
#    "_GoBack" bookmark that used to sit next to "Frankie Homewood" at the
#    end of the letter, because a document may only contain one bookmark of
#    a given name - adding a new one automatically removes the old one.
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$runStart = $full.IndexOf(" skills. Gained through my extensive experience performing in a campanology (tower bell ringing) group")

$fragments = @(
    " skills. Gained through",
    " my extensive experience performing in a campanology (",
    "tower",
    " ",
    "b",
    "ell ",
    "r",
    "ing",
    "ing",
    ") group"
)

$offset = 0
for ($i = 0; $i -lt $fragments.Length - 1; $i++) {
    $offset += $fragments[$i].Length
    $pos = $runStart + $offset
    $rSplit = $d.Range($pos, $pos)
    if ($i -eq 0) {
        # This is the boundary right after "through" - put the real bookmark here.
        $d.Bookmarks.Add("_GoBack", $rSplit) | Out-Null
    } else {
        $d.Bookmarks.Add("ZZZ_split" + $i, $rSplit) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 4. Every fragment inherited the xml:space="preserve" flag of the single
#    run it was split from, even the ones with no leading/trailing space
#    ("tower", "b", "r", "ing", "ing", ") group"). Force Word to recompute
#    that flag for just those fragments by writing different text into them
#    and then writing the correct text straight back - a genuine content
#    change on an already isolated range recalculates the flag correctly
#    without touching any neighbouring run. Positions are tracked by walking
#    through the fragment list in order (rather than re-searching for short,
#    ambiguous substrings such as "b" or "ing") so the exact right fragment
#    is targeted every time.
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$runStart = $full.IndexOf(" skills. Gained through my extensive experience performing in a campanology (tower bell ringing) group")

$offset = 0
for ($i = 0; $i -lt $fragments.Length; $i++) {
    $frag = $fragments[$i]
    $fragStart = $runStart + $offset
    if ($frag -ne $frag.Trim() -or $frag -eq "") {
        # Fragment genuinely needs (or already correctly has) xml:space="preserve";
        # leave it untouched.
        $offset += $frag.Length
        continue
    }
    $rFrag = $d.Range($fragStart, $fragStart + $frag.Length)
    $rFrag.Text = "ZZZTMPZZZ"
    $rFrag2 = $d.Range($fragStart, $fragStart + 9)
    $rFrag2.Text = $frag
    $offset += $frag.Length
}

# ---------------------------------------------------------------------------
# 5. Remove all temporary bookmarks used purely to force the run splits.
#    ("_GoBack" is intentionally left in the document.)
# ---------------------------------------------------------------------------

$tempNames = @("ZZZ_splitBefore", "ZZZ_splitAfter")
for ($i = 1; $i -lt $fragments.Length - 1; $i++) {
    $tempNames += ("ZZZ_split" + $i)
}
foreach ($name in $tempNames) {
    if ($d.Bookmarks.Exists($name)) {
        $d.Bookmarks.Item($name).Delete()
    }
}

Write-Output $d.Content.Text
